$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("N3").Value = 1.65
$ws.Range("O3").Value = 2.2

# Row 9 updates
$ws.Range("G9").Value = 5.4
$ws.Range("H9").Value = 4
$ws.Range("M9").Value = 3.65
$ws.Range("R9").Value = 1.7
$ws.Range("S9").Value = 1.93
$ws.Range("T9").Value = 17.5
$ws.Range("V9").Value = 17
$ws.Range("W9").Value = 100
$ws.Range("X9").Value = 50
$ws.Range("Y9").Value = 45
$ws.Range("Z9").Value = 13
$ws.Range("AA9").Value = 8
$ws.Range("AD9").Value = 400
$ws.Range("AE9").Value = 7.8
$ws.Range("AF9").Value = 7.8
$ws.Range("AG9").Value = 8
$ws.Range("AH9").Value = 11.25
$ws.Range("AI9").Value = 11.75
